$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append new row 20 to the "Logs" sheet ---
$ws.Range("A20").Value = "Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$ws.Range("B20").Value = "mailmind.test@zohomail.eu"
$ws.Range("C20").Value = "Testmail #18: Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$ws.Range("D20").Value = "Bestelling / Levering"

$antwoord20 = @"
Beste klant,
Bedankt voor je e-mail. Helaas kan ik je in dit geval niet verder helpen met je bestelling van 200 stuks M8-bouten RVS voor Van Dijk. Ik raad je aan om contact op te nemen met onze verkoopafdeling of een van onze vertegenwoordigers, zodat zij je verder kunnen assisteren met het plaatsen van deze bestelling.
Mocht je nog andere vragen hebben of hulp nodig hebben, laat het ons gerust weten.
Met vriendelijke groet,
[Naam] E-mailassistent bij [Bedrijfsnaam]
"@

$ws.Range("E20").Value = $antwoord20
$ws.Range("F20").Value = "2025-07-31 22:01:55"
$ws.Range("G20").Value = "Ja"
$ws.Range("H20").Value = "Nee"
$ws.Range("I20").Value = "Ja"
$ws.Range("J20").Value = "Nee"

# --- Extend the conditional formatting ranges from row 19 to row 20 ---
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = "$($col)2:$($col)19"
    $newRange = "$($col)2:$($col)20"
    $fcs = $ws.Range($oldRange).FormatConditions
    $cnt = $fcs.Count()
    for ($i = 1; $i -le $cnt; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($ws.Range($newRange))
    }
}

# --- Update the "Dashboard" summary: swap Bestelling / Retour rows ---
$dash.Range("A5").Value = "Bestelling / Levering"
$dash.Range("B5").Value = 2
$dash.Range("A6").Value = "Retour / Terugbetaling"
$dash.Range("B6").Value = 2
